# Apply updated odds values to Sheet1, matching the target diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 2 updates
$ws.Range("G2").Value = 2.15
$ws.Range("H2").Value = 2.57
$ws.Range("I2").Value = 4.4
$ws.Range("J2").Value = 2.92
$ws.Range("L2").Value = 5.1
$ws.Range("N2").Value = 4.3
$ws.Range("S2").Value = 5.7
$ws.Range("Y2").Value = 4.7
$ws.Range("AB2").Value = 21
$ws.Range("AE2").Value = 4.3
$ws.Range("AF2").Value = 5.5
$ws.Range("AI2").Value = 7.9
$ws.Range("AJ2").Value = 22
$ws.Range("AK2").Value = 16
$ws.Range("AL2").Value = 90
$ws.Range("AM2").Value = 65

# Row 3 updates
$ws.Range("G3").Value = 1.93
$ws.Range("H3").Value = 3.7
$ws.Range("I3").Value = 3.3
$ws.Range("J3").Value = 2.45
$ws.Range("L3").Value = 3.65
$ws.Range("S3").Value = 2.5
$ws.Range("T3").Value = 1.4
$ws.Range("Y3").Value = 9
$ws.Range("Z3").Value = 10.25
$ws.Range("AB3").Value = 17
$ws.Range("AC3").Value = 14.5
$ws.Range("AE3").Value = 13
$ws.Range("AF3").Value = 7.4
$ws.Range("AI3").Value = 12
$ws.Range("AJ3").Value = 19
$ws.Range("AK3").Value = 11.5
$ws.Range("AM3").Value = 26
$ws.Range("AN3").Value = 30

# Row 4 updates
$ws.Range("G4").Value = 1.65
$ws.Range("H4").Value = 4
$ws.Range("L4").Value = 4.75
$ws.Range("M4").Value = 1.03
$ws.Range("N4").Value = 10
$ws.Range("Q4").Value = 1.67
$ws.Range("R4").Value = 2.15
$ws.Range("U4").Value = 1.3
$ws.Range("Y4").Value = 8.5
$ws.Range("Z4").Value = 8.5
$ws.Range("AF4").Value = 8
$ws.Range("AJ4").Value = 26

# Row 5 updates
$ws.Range("S5").Value = 2
$ws.Range("T5").Value = 1.73
$ws.Range("U5").Value = 1.19

# Row 6 updates
$ws.Range("G6").Value = 1.27
$ws.Range("H6").Value = 6
$ws.Range("I6").Value = 9
$ws.Range("L6").Value = 7.5
$ws.Range("Q6").Value = 1.44
$ws.Range("R6").Value = 2.63
$ws.Range("U6").Value = 1.25
$ws.Range("Y6").Value = 9
$ws.Range("Z6").Value = 7
$ws.Range("AA6").Value = 9
$ws.Range("AF6").Value = 12
$ws.Range("AI6").Value = 26
$ws.Range("AK6").Value = 26
$ws.Range("AL6").Value = 101

# Row 7 updates
$ws.Range("H7").Value = 3.3
$ws.Range("I7").Value = 2.92
$ws.Range("L7").Value = 3.45
$ws.Range("N7").Value = 6.8
$ws.Range("T7").Value = 1.26
$ws.Range("X7").Value = 1.88
$ws.Range("AB7").Value = 23
$ws.Range("AE7").Value = 6.8
$ws.Range("AF7").Value = 6.3
$ws.Range("AK7").Value = 10.75
$ws.Range("AM7").Value = 26
$ws.Range("AO7").Value = 600
